$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Plain text / non-ambiguous-number cells: Coin name, Link URL, Volume% text
$textUpdates = @{
    "E2" = "  +4.03%  "
    "E3" = "  +3.81%  "
    "E4" = "  -0.09%  "
    "E5" = "  +1.98%  "
    "E6" = "  -0.10%  "
    "E7" = "  -0.98%  "
    "E8" = "  +1.18%  "
    "E9" = "  +1.57%  "
    "E10" = "  +3.44%  "
    "E11" = "  +0.10%  "
    "E12" = "  +4.49%  "
    "E13" = "  +0.82%  "
    "E14" = "  +0.74%  "
    "E15" = "  +2.82%  "
    "E16" = "  -0.09%  "
    "E17" = "  -0.11%  "
    "E18" = "  +4.08%  "
    "E19" = "  +1.08%  "
    "E20" = "  +1.85%  "
    "E21" = "  +3.41%  "
    "E22" = "  +2.91%  "
    "E23" = "  +1.86%  "
    "E24" = "  +0.37%  "
    "E25" = "  +1.14%  "
    "E26" = "  +1.48%  "
    "E27" = "  +0.45%  "
    "E28" = "  +4.74%  "
    "E29" = "  +2.72%  "
    "E30" = "  +1.50%  "
    "E31" = "  +2.34%  "
    "E32" = "  +1.55%  "
    "E33" = "  +2.31%  "
    "E34" = "  +0.31%  "
    "E35" = "  +3.50%  "
    "E36" = "  +1.59%  "
    "E37" = "  +18.65%  "
    "E38" = "  +7.90%  "
    "E39" = "  -7.01%  "
    "B40" = "PaxDollar"
    "C40" = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
    "E40" = "  +0.05%  "
    "B41" = "RenderToken"
    "C41" = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
    "E41" = "  +4.39%  "
    "E42" = "  -1.77%  "
    "E43" = "  +15.48%  "
    "E44" = "  +1.80%  "
    "E45" = "  +4.76%  "
    "E46" = "  +1.48%  "
    "E47" = "  +0.86%  "
    "E48" = "  +4.18%  "
    "E49" = "  +4.49%  "
    "E50" = "  +1.93%  "
    "E51" = "  +2.09%  "
}

# Price cells that look numeric (e.g. "239.90", "5.340") -- must be forced to
# text so Excel keeps the exact string (trailing zeros, dotted-thousands, etc.)
# instead of silently coercing them to a Double.
$priceUpdates = @{
    "D2" = "26.348.36"
    "D3" = "1.722.95"
    "D5" = "239.90"
    "D6" = "0.9997"
    "D7" = "0.4730"
    "D8" = "0.2632"
    "D9" = "0.06239"
    "D10" = "1.716.82"
    "D11" = "0.07076"
    "D12" = "15.32"
    "D13" = "0.5944"
    "D14" = "4.411"
    "D15" = "76.37"
    "D18" = "26.344.56"
    "D19" = "0.000006812"
    "D20" = "11.61"
    "D21" = "1.936.48"
    "D22" = "4.563"
    "D23" = "8.788"
    "D24" = "5.340"
    "D25" = "135.04"
    "D26" = "15.26"
    "D28" = "1.765"
    "D29" = "106.88"
    "D30" = "4.022"
    "D31" = "3.704"
    "D32" = "0.07758"
    "D33" = "0.04459"
    "D34" = "2.611"
    "D35" = "0.9768"
    "D36" = "0.6218"
    "D37" = "115.98"
    "D38" = "0.9248"
    "D39" = "2.424"
    "D40" = "1.000"
    "D41" = "1.907"
    "D42" = "0.01472"
    "D43" = "5.358"
    "D44" = "0.3828"
    "D46" = "6.270"
    "D47" = "0.05291"
    "D48" = "30.61"
    "D49" = "7.667"
    "D50" = "0.3398"
    "D51" = "1.221"
}

foreach ($addr in $textUpdates.Keys) {
    $ws.Range($addr).Value = $textUpdates[$addr]
}

foreach ($addr in $priceUpdates.Keys) {
    $ws.Range($addr).NumberFormat = "@"
    $ws.Range($addr).Value = $priceUpdates[$addr]
}
